$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$question = "HTML stands for -"
$optA = "HighText Machine Language"
$optB = "HyperText and links Markup Language"
$optC = "HyperText Markup Language"
$optD = "None of these"
$answer = "HyperText Markup Language"

for ($r = 1; $r -le 5; $r++) {
    $ws.Cells.Item($r, 1).Value = $question
    $ws.Cells.Item($r, 2).Value = $optA
    $ws.Cells.Item($r, 3).Value = $optB
    $ws.Cells.Item($r, 4).Value = $optC
    $ws.Cells.Item($r, 5).Value = $optD
    $ws.Cells.Item($r, 6).Value = $answer
}

$ws.Range("F8").Select()

$ws.PageSetup.Orientation = 1
